{"js": "// The document contains paragraphs whose `<id>...</id>` marker was\n// previously split into three separate runs:\n//   run0 \"<id>\"        (Courier New, color 7f6000, sz 18 - the \"tag\" style)\n//   run1 \"p035v_N\"      (default font, color 000000 - the \"value\" style)\n//   run2 \"</id>\"        (Courier New, color 7f6000, sz 18 - the \"tag\" style)\n// The edit merges those three runs into a single run per paragraph whose\n// text is the concatenation \"<id>p035v_N</id>\" and whose formatting is the\n// formatting of the first (\"<id>\") run. Only the paragraphs identifying the\n// transcription entries themselves (p035v_1 / p035v_2) are affected; the\n// figure id paragraph (fig_p035v_1) is left untouched.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst targets = [\"<id>p035v_1</id>\", \"<id>p035v_2</id>\"];\n\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const paragraph = paragraphs.items[i];\n  if (targets.indexOf(paragraph.text) !== -1) {\n    // Replacing the paragraph's whole range with the same literal text\n    // collapses the existing runs into one run that carries the\n    // formatting of the first original run (matching the target XML).\n    paragraph.insertText(paragraph.text, \"Replace\");\n  }\n}\n\nawait context.sync();\n", "ps1": "# The `<id>...</id>` marker in two transcription-entry paragraphs was\n# previously split across three separate runs:\n#   \"<id>\"        (Courier New, color 7f6000, sz 18 - the \"tag\" style)\n#   \"p035v_N\"      (default font, color 000000 - the \"value\" style)\n#   \"</id>\"        (Courier New, color 7f6000, sz 18 - the \"tag\" style)\n# Using Find & Replace (ReplaceAll) against the exact literal text merges\n# the matched span back into a single run that carries the formatting of\n# the first original run, matching the target XML. The figure-id\n# paragraph (\"<id>fig_p035v_1</id>\") is a different literal string, so it\n# is left untouched.\n\n$d = $word.ActiveDocument\n\n$ids = @(\"p035v_1\", \"p035v_2\")\n\nforeach ($id in $ids) {\n    $needle = \"<id>\" + $id + \"</id>\"\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Execute($needle, $false, $false, $false, $false, $false, $true, 1, $false, $needle, 2)\n}\n"}
